$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The financial-statement table (period headers in row 8, publish dates in row
# 9, and the data rows 11-27) is quarterly/cumulative data that rolls forward
# by one reporting period: the oldest period (column D) is dropped, every
# later period shifts one column to the left, and a brand-new latest period
# is appended in column M.
#
# Do the shift with a native copy/paste of values+formats (not a per-cell
# literal re-type) so every cell keeps its original style/number-format,
# and so date-look-alike text (e.g. "1401-10-28") is carried over as the
# plain text it already is instead of being re-interpreted by value-entry
# parsing.
$ws.Range("E8:M27").Copy()
$ws.Range("D8:L27").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# --- New, latest reporting period: column M ---

# Row 8: period-end label; Row 9: publish-date label for that period.
$ws.Cells.Item(8, 13).Value = "12 ماهه منتهی به 1401/12"
$ws.Cells.Item(9, 13).Value = "1402-02-06 (2)"

# Rows 11-27: the new period's reported financial figures.
$ws.Cells.Item(11, 13).Value = 5516295
$ws.Cells.Item(12, 13).Value = -2070101
$ws.Cells.Item(13, 13).Value = 3446194
$ws.Cells.Item(14, 13).Value = -231795
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(16, 13).Value = -15654
$ws.Cells.Item(17, 13).Value = 3198745
$ws.Cells.Item(18, 13).Value = -62657
$ws.Cells.Item(19, 13).Value = 84016
$ws.Cells.Item(20, 13).Value = 3220104
$ws.Cells.Item(21, 13).Value = -458085
$ws.Cells.Item(22, 13).Value = 2762019
$ws.Cells.Item(23, 13).Value = 0
$ws.Cells.Item(24, 13).Value = 2762019
$ws.Cells.Item(25, 13).Value = 9207
$ws.Cells.Item(26, 13).Value = 300000
$ws.Cells.Item(27, 13).Value = 9207

# Row 9 also carries one mid-table correction beyond the pure shift: what had
# been column J's publish-date annotation changes from "1401-10-28 (7)" to
# "1402-02-06 (9)" once it lands in column I after the shift.
$ws.Cells.Item(9, 9).Value = "1402-02-06 (9)"
